# Update cryptos list with latest prices and 1h volume changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'42.098.11"
$ws.Range('E2').Value = '  -0.84%  '
$ws.Range('D3').Value = "'2.224.30"
$ws.Range('E3').Value = '  -1.03%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = "'242.82"
$ws.Range('E5').Value = '  -1.08%  '
$ws.Range('D6').Value = "'0.627"
$ws.Range('E6').Value = '  -0.24%  '
$ws.Range('D7').Value = "'74.53"
$ws.Range('E7').Value = '  -1.92%  '
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('D9').Value = "'0.605"
$ws.Range('E9').Value = '  -3.07%  '
$ws.Range('D10').Value = "'42.63"
$ws.Range('E10').Value = '  -4.18%  '
$ws.Range('D11').Value = "'0.0958"
$ws.Range('E11').Value = '  +0.36%  '
$ws.Range('D12').Value = "'6.96"
$ws.Range('E12').Value = '  -4.00%  '
$ws.Range('E13').Value = '  +0.15%  '
$ws.Range('D14').Value = "'2.561.59"
$ws.Range('E14').Value = '  -1.04%  '
$ws.Range('D15').Value = "'14.29"
$ws.Range('E15').Value = '  -2.42%  '
$ws.Range('D16').Value = "'0.837"
$ws.Range('E16').Value = '  -2.90%  '
$ws.Range('D17').Value = "'2.218.04"
$ws.Range('E17').Value = '  -1.76%  '
$ws.Range('D18').Value = "'41.993.68"
$ws.Range('E18').Value = '  -0.85%  '
$ws.Range('E19').Value = '  +3.64%  '
$ws.Range('D20').Value = "'6.21"
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').Value = "'72.95"
$ws.Range('E21').Value = '  +0.89%  '
$ws.Range('D22').Value = "'11.15"
$ws.Range('E22').Value = '  -2.71%  '
$ws.Range('D23').Value = "'230.29"
$ws.Range('D24').Value = "'2.09"
$ws.Range('E24').Value = '  -7.49%  '
$ws.Range('D25').Value = "'1.00"
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range('D26').Value = "'11.40"
$ws.Range('E26').Value = '  -3.83%  '
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('E28').Value = '  -2.04%  '
$ws.Range('E29').Value = '  -3.35%  '
$ws.Range('D30').Value = "'166.49"
$ws.Range('E30').Value = '  -0.74%  '
$ws.Range('D31').Value = "'20.57"
$ws.Range('E31').Value = '  -0.88%  '
$ws.Range('D32').Value = "'5.65"
$ws.Range('E32').Value = '  -5.09%  '
$ws.Range('D33').Value = "'0.0800"
$ws.Range('E33').Value = '  -2.78%  '
$ws.Range('D34').Value = "'30.05"
$ws.Range('E34').Value = '  -2.92%  '
$ws.Range('E35').Value = '  -1.23%  '
$ws.Range('D36').Value = "'0.109"
$ws.Range('E36').Value = '  -8.45%  '
$ws.Range('D37').Value = "'4.32"
$ws.Range('E37').Value = '  -8.74%  '
$ws.Range('D38').Value = "'0.0303"
$ws.Range('E38').Value = '  -3.79%  '
$ws.Range('D39').Value = "'13.18"
$ws.Range('E39').Value = '  -5.88%  '
$ws.Range('D40').Value = "'2.13"
$ws.Range('E40').Value = '  -2.86%  '
$ws.Range('B41').Value = 'MultiversX'
$ws.Range('C41').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D41').Value = "'65.19"
$ws.Range('E41').Value = '  +1.64%  '
$ws.Range('B42').Value = 'THORChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D42').Value = "'5.69"
$ws.Range('E42').Value = '  -1.74%  '
$ws.Range('D43').Value = "'0.198"
$ws.Range('E43').Value = '  -2.34%  '
$ws.Range('D44').Value = "'8.71"
$ws.Range('E44').Value = '  -2.02%  '
$ws.Range('D45').Value = "'104.13"
$ws.Range('E45').Value = '  -3.41%  '
$ws.Range('E46').Value = '  -2.07%  '
$ws.Range('D47').Value = "'2.35"
$ws.Range('E47').Value = '  -3.91%  '
$ws.Range('E48').Value = '  -2.69%  '
$ws.Range('E49').Value = '  -1.89%  '
$ws.Range('E50').Value = '  -1.07%  '
$ws.Range('D51').Value = "'2.435.62"
$ws.Range('E51').Value = '  -1.17%  '
